$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "date" placeholder on the slide master + every slide layout
#    cached field text: 22.03.2012 -> 23.03.2012
# ---------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "22.03.2012") {
                $shp.TextFrame.TextRange.Text = "23.03.2012"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 6, textbox "Textfeld 117": widen it and reword
#    "Mapping depends on enclosed elements" ->
#    "Mapping depending on enclosed elements" (with " on " split into
#    two runs: " " and "on ")
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$shp = $slide6.Shapes("Textfeld 117")

$tr = $shp.TextFrame.TextRange
$dependsRange = $tr.Characters(9, 7)   # "depends"
$dependsRange.Text = "depending"

$onRange = $tr.Characters(19, 3)       # "on " (tail of the former " on " run)
$onRange.Text = "on "

$shp.Width = 106.17686
